$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1738.4445
$ws.Range("J40").Value = 1472.4375
$ws.Range("L40").Value = 1472.4375
$ws.Range("N40").Value = -1822.4375
$ws.Range("H62").Value = 3081.1667
$ws.Range("I62").Value = 3668.3333
$ws.Range("K62").Value = 3668.3333
$ws.Range("M62").Value = -3044.3333
$ws.Range("H65").Value = 3081.1667
$ws.Range("I65").Value = 3668.3333
$ws.Range("K65").Value = 18341.6665
$ws.Range("M65").Value = -15221.6665
$ws.Range("H121").Value = 710.97144
$ws.Range("J121").Value = 727.625
$ws.Range("L121").Value = 2182.875
$ws.Range("N121").Value = -5676.875
$ws.Range("H123").Value = 15150
$ws.Range("J123").Value = 15150
$ws.Range("L123").Value = 15150
$ws.Range("N123").Value = -24950
$ws.Range("H137").Value = 21278152
$ws.Range("I137").Value = 1207.963
$ws.Range("K137").Value = 3623.889
$ws.Range("M137").Value = -1073.889
$ws.Range("H138").Value = 2656.21
$ws.Range("I138").Value = 1277.619
$ws.Range("J138").Value = 3022.671
$ws.Range("K138").Value = 3832.857
$ws.Range("L138").Value = 9068.012999999999
$ws.Range("M138").Value = 1307.143
$ws.Range("N138").Value = -19348.013
$ws.Range("H140").Value = 34000
$ws.Range("J140").Value = 34000
$ws.Range("L140").Value = 34000
$ws.Range("N140").Value = -44360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25840.527
$ws.Range("I32").Value = 24282.516
$ws.Range("J32").Value = 33890.25
$ws.Range("K32").Value = 24282.516
$ws.Range("L32").Value = 33890.25
$ws.Range("M32").Value = -23995.516
$ws.Range("N32").Value = -34464.25
$ws.Range("H122").Value = 1586.6666
$ws.Range("I122").Value = 1344
$ws.Range("K122").Value = 4032
$ws.Range("M122").Value = -1582
$ws.Range("H132").Value = 7676.1025
$ws.Range("I132").Value = 9180.821
$ws.Range("J132").Value = 3845.9092
$ws.Range("K132").Value = 27542.463
$ws.Range("L132").Value = 11537.7276
$ws.Range("M132").Value = -25012.463
$ws.Range("N132").Value = -16597.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1500
$ws.Range("I17").Value = 1500
$ws.Range("K17").Value = 1500
$ws.Range("M17").Value = -1326
$ws.Range("H31").Value = 37044136
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 37044136
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H132").Value = 3572949.5
$ws.Range("I132").Value = 1307.8
$ws.Range("J132").Value = 25002800
$ws.Range("K132").Value = 3923.4
$ws.Range("L132").Value = 75008400
$ws.Range("M132").Value = -1393.4
$ws.Range("N132").Value = -75013460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 322.4074
$ws.Range("I2").Value = 359.58334
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 2157.50004
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -2044.50004
$ws.Range("N2").Value = -376
$ws.Range("H17").Value = 201
$ws.Range("I17").Value = 201
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 603
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -434
$ws.Range("N17").ClearContents()
$ws.Range("H34").Value = 1158.8334
$ws.Range("I34").Value = 225
$ws.Range("J34").Value = 1625.75
$ws.Range("K34").Value = 675
$ws.Range("L34").Value = 4877.25
$ws.Range("M34").Value = -591
$ws.Range("N34").Value = -5045.25
$ws.Range("H39").Value = 3775
$ws.Range("J39").Value = 3775
$ws.Range("L39").Value = 11325
$ws.Range("N39").Value = -11913
$ws.Range("H55").Value = 2508
$ws.Range("J55").Value = 2508
$ws.Range("L55").Value = 7524
$ws.Range("N55").Value = -7878
$ws.Range("H131").Value = 784.36
$ws.Range("I131").Value = 482.5
$ws.Range("J131").Value = 810.6087
$ws.Range("K131").Value = 1447.5
$ws.Range("L131").Value = 2431.8261
$ws.Range("M131").Value = 3592.5
$ws.Range("N131").Value = -12511.8261
$ws.Range("H132").Value = 1720.0869
$ws.Range("I132").Value = 888.6667
$ws.Range("J132").Value = 2627.0908
$ws.Range("K132").Value = 7998.0003
$ws.Range("L132").Value = 23643.8172
$ws.Range("M132").Value = -5468.0003
$ws.Range("N132").Value = -28703.8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 20002216
$ws.Range("I122").Value = 35716490
$ws.Range("J122").Value = 2233.1365
$ws.Range("K122").Value = 107149470
$ws.Range("L122").Value = 6699.4095
$ws.Range("M122").Value = -107147020
$ws.Range("N122").Value = -11599.4095
$ws.Range("H126").Value = 7530
$ws.Range("I126").Value = 12800
$ws.Range("K126").Value = 38400
$ws.Range("M126").Value = -35930
$ws.Range("H132").Value = 5849.467
$ws.Range("I132").Value = 6727
$ws.Range("J132").Value = 3436.25
$ws.Range("K132").Value = 20181
$ws.Range("L132").Value = 10308.75
$ws.Range("M132").Value = -17651
$ws.Range("N132").Value = -15368.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1384.5217
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 1440.1904
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 1440.1904
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1816.1904
$ws.Range("H61").Value = 19609546
$ws.Range("I61").Value = 1626.375
$ws.Range("J61").Value = 37038810
$ws.Range("K61").Value = 1626.375
$ws.Range("L61").Value = 37038810
$ws.Range("M61").Value = -1424.375
$ws.Range("N61").Value = -37039214
$ws.Range("H113").Value = 19609546
$ws.Range("I113").Value = 1626.375
$ws.Range("J113").Value = 37038810
$ws.Range("K113").Value = 1626.375
$ws.Range("L113").Value = 37038810
$ws.Range("M113").Value = 543.625
$ws.Range("N113").Value = -37043150
$ws.Range("H122").Value = 10300.8
$ws.Range("I122").Value = 14668
$ws.Range("K122").Value = 44004
$ws.Range("M122").Value = -41554
$ws.Range("H136").Value = 7279.391
$ws.Range("I136").Value = 8180.316
$ws.Range("K136").Value = 24540.948
$ws.Range("M136").Value = -21990.948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 23249
$ws.Range("J63").Value = 23249
$ws.Range("L63").Value = 23249
$ws.Range("N63").Value = -24497
$ws.Range("H66").Value = 23249
$ws.Range("J66").Value = 23249
$ws.Range("L66").Value = 69747
$ws.Range("N66").Value = -75987
$ws.Range("H113").Value = 444.75
$ws.Range("I113").Value = 449.25
$ws.Range("J113").Value = 440.25
$ws.Range("K113").Value = 1347.75
$ws.Range("L113").Value = 1320.75
$ws.Range("M113").Value = 822.25
$ws.Range("N113").Value = -5660.75
$ws.Range("H126").Value = 83346424
$ws.Range("I126").Value = 125017064
$ws.Range("J126").Value = 5141
$ws.Range("K126").Value = 375051192
$ws.Range("L126").Value = 15423
$ws.Range("M126").Value = -375048722
$ws.Range("N126").Value = -20363
$ws.Range("H136").Value = 1656.5227
$ws.Range("I136").Value = 1557.225
$ws.Range("J136").Value = 2649.5
$ws.Range("K136").Value = 4671.674999999999
$ws.Range("L136").Value = 7948.5
$ws.Range("M136").Value = -2121.674999999999
$ws.Range("N136").Value = -13048.5
